# Apply updated cryptocurrency price/volume data to the "cryptos" worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. "D" column values are price strings that
# often look numeric (e.g. "336.66"), so they are written through a small helper that
# forces them to remain plain text (matching the original inline-string cells) instead
# of being auto-converted into numbers by Excel.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" '30.587.27'
$ws.Range("E2").Value = '  +0.44%  '
Set-TextValue "D3" '2.112.23'
$ws.Range("E3").Value = '  +0.82%  '
$ws.Range("E4").Value = '  +1.02%  '
Set-TextValue "D5" '336.66'
$ws.Range("E5").Value = '  +2.18%  '
Set-TextValue "D6" '1.011'
$ws.Range("E6").Value = '  +0.93%  '
Set-TextValue "D7" '0.5252'
$ws.Range("E7").Value = '  +0.81%  '
Set-TextValue "D8" '0.4556'
$ws.Range("E8").Value = '  +4.88%  '
Set-TextValue "D9" '55.14'
$ws.Range("E9").Value = '  +5.43%  '
Set-TextValue "D10" '0.09017'
$ws.Range("E10").Value = '  +1.93%  '
Set-TextValue "D11" '1.169'
$ws.Range("E11").Value = '  +1.15%  '
Set-TextValue "D12" '24.51'
$ws.Range("E12").Value = '  +0.50%  '
Set-TextValue "D13" '2.117.15'
$ws.Range("E13").Value = '  +1.42%  '
Set-TextValue "D14" '6.869'
$ws.Range("E14").Value = '  +2.78%  '
Set-TextValue "D15" '8.129'
$ws.Range("E15").Value = '  +5.85%  '
Set-TextValue "D16" '0.00001177'
$ws.Range("E16").Value = '  +5.17%  '
Set-TextValue "D17" '97.17'
$ws.Range("E17").Value = '  +1.33%  '
Set-TextValue "D18" '1.011'
$ws.Range("E18").Value = '  +0.89%  '
Set-TextValue "D19" '0.06690'
Set-TextValue "D20" '19.39'
$ws.Range("E20").Value = '  +0.96%  '
$ws.Range("E21").Value = '  +0.88%  '
Set-TextValue "D22" '6.262'
$ws.Range("E22").Value = '  -0.03%  '
Set-TextValue "D23" '30.648.28'
$ws.Range("E23").Value = '  +0.59%  '
Set-TextValue "D24" '12.86'
$ws.Range("E24").Value = '  +5.62%  '
Set-TextValue "D25" '2.361'
$ws.Range("E25").Value = '  +1.18%  '
Set-TextValue "D26" '2.358.26'
$ws.Range("E26").Value = '  +1.11%  '
Set-TextValue "D27" '22.34'
$ws.Range("E27").Value = '  +0.36%  '
Set-TextValue "D28" '163.54'
$ws.Range("E28").Value = '  +0.69%  '
Set-TextValue "D29" '2.529'
$ws.Range("E29").Value = '  -2.43%  '
Set-TextValue "D30" '133.74'
$ws.Range("E30").Value = '  +1.74%  '
Set-TextValue "D31" '1.225'
$ws.Range("E31").Value = '  +2.82%  '
Set-TextValue "D32" '0.1073'
$ws.Range("E32").Value = '  +0.43%  '
Set-TextValue "D33" '6.353'
$ws.Range("E33").Value = '  +3.45%  '
Set-TextValue "D34" '1.628'
$ws.Range("E34").Value = '  -2.12%  '
Set-TextValue "D35" '3.969'
$ws.Range("E35").Value = '  +2.18%  '
$ws.Range("E36").Value = '  +3.48%  '
Set-TextValue "D37" '5.883'
$ws.Range("E37").Value = '  +7.81%  '
Set-TextValue "D38" '0.02601'
$ws.Range("E38").Value = '  +1.33%  '
Set-TextValue "D39" '0.06829'
$ws.Range("E39").Value = '  +0.44%  '
Set-TextValue "D40" '0.2328'
$ws.Range("E40").Value = '  +2.85%  '
Set-TextValue "D41" '12.62'
$ws.Range("E41").Value = '  -0.85%  '
Set-TextValue "D42" '0.6859'
$ws.Range("E42").Value = '  -0.74%  '
Set-TextValue "D43" '1.254'
$ws.Range("E43").Value = '  -0.91%  '
Set-TextValue "D44" '0.6449'
$ws.Range("E44").Value = '  +1.07%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D45" '14.13'
$ws.Range("E45").Value = '  +1.43%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D46" '2.315'
$ws.Range("E46").Value = '  +5.16%  '
Set-TextValue "D47" '3.680'
$ws.Range("E47").Value = '  +1.55%  '
Set-TextValue "D48" '1.252'
$ws.Range("E48").Value = '  +0.87%  '
Set-TextValue "D49" '0.00000000347'
$ws.Range("E49").Value = '  +18.29%  '
$ws.Range("E50").Value = '  -1.71%  '
Set-TextValue "D51" '83.28'
$ws.Range("E51").Value = '  +1.73%  '
